# Edit script: "Se agregan cambios a la semana"
# 1. Merge the split "<Day> " + "<N>" runs into a single run for the
#    "Martes 21", "Miercoles 22" and "Jueves 23" red date headers.
# 2. Append six new paragraphs (three red date headers and three list
#    entries describing what happened those days) right before the
#    trailing empty paragraph at the end of the document.

$d = $word.ActiveDocument

# --- Part 1: merge the two-run date headers into a single run -------------
$mergeTargets = @("Martes 21", "Mi" + [char]0xE9 + "rcoles 22", "Jueves 23")
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs($i)
    $paraText = $para.Range.Text.TrimEnd([char]13)
    foreach ($target in $mergeTargets) {
        if ($paraText -eq $target) {
            $para.Range.Find.Execute($target, $false, $false, $false, $false, `
                $false, $true, 1, $false, $target, 2) | Out-Null
        }
    }
}

# --- Part 2: append the new week-end paragraphs ----------------------------
# The document always ends with one empty paragraph (w:ind w:left="360").
# Collapse a range to the very start of that paragraph and InsertXML the
# new paragraphs there, so the trailing empty paragraph is left intact
# right after them.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$insertionPoint = $lastPara.Range
$insertionPoint.Collapse(1)

$newParagraphsXml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>Viernes 24</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>Sin cambios</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>Sábado 25 y domingo 26</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>Se cumplen los horarios establecidos, finalmente se agrega mas tiempo en la noche de 8pm a 10pm para repasar métodos numéricos en grupo de estudio, para el parcial del día siguiente.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>Lunes 27</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>Día de parcial, luego del trabajo descanso ya que no hay temas para repasar sobre la asignatura de métodos numéricos y me encuentro bastante cansado.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertionPoint.InsertXML($newParagraphsXml)

Write-Output "done"
